$d = $word.ActiveDocument

$bodyXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Circle Language Spec Plan</w:t></w:r><w:r><w:br/></w:r><w:r><w:t xml:space="preserve">Programming </w:t></w:r><w:r><w:t>Language Programmed Within Itself</w:t></w:r><w:r><w:t xml:space="preserve"> Spec</w:t></w:r><w:r><w:br/><w:t>2008-03</w:t></w:r><w:r><w:br/></w:r><w:r><w:t xml:space="preserve">Project </w:t></w:r><w:r><w:t>Summary</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:ind w:left="142"/><w:rPr><w:i/><w:iCs/><w:sz w:val="20"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="20"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">Author: </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="20"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">JJ </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="20"/><w:szCs w:val="22"/></w:rPr><w:t>van Zon</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="142"/><w:rPr><w:i/><w:iCs/><w:sz w:val="20"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="20"/><w:szCs w:val="22"/></w:rPr><w:t>Location: Oosterhout</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="20"/><w:szCs w:val="22"/></w:rPr><w:t>, The Netherlands</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="SpacingCharChar"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:r><w:t>Goal</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="SpacingCharChar"/><w:ind w:left="284"/><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:t>Write</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> the article </w:t></w:r><w:r><w:rPr><w:i/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Programming </w:t></w:r><w:r><w:rPr><w:i/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:t>Language Programmed Within Itself</w:t></w:r><w:r><w:rPr><w:i/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:t>which produced version  2008-03-10 00  1.0 .</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:r><w:t>Super-</w:t></w:r><w:r><w:t>P</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>roject</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="284"/><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr><w:bookmarkStart w:id="1" w:name="_Hlk37947535"/><w:r><w:t>This is a sub-project of the super-project</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Circle Language Spec, Fundamental Principles Spec Part A</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:bookmarkEnd w:id="1"/><w:r><w:t>Date</w:t></w:r><w:r><w:t xml:space="preserve"> and Time</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="284"/></w:pPr><w:r><w:t>2008-03-06 &#8211; 2008-03-09</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="284"/></w:pPr><w:r><w:t>8 hours of work</w:t></w:r></w:p><w:p/><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$d.Content.InsertXML($bodyXml)

$h3 = $d.Styles("Heading3")
$h3.Font.Italic = 0
$h3.Font.Size = 10
$h3.ParagraphFormat.SpaceBefore = 9
$h3.ParagraphFormat.SpaceAfter = 9

Write-Output "done"
